# Apply the latest cryptos snapshot update (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.560.02'
$ws.Range("E2").Value = '  +0.09%  '

$ws.Range("D3").Value = '1.914.47'
$ws.Range("E3").Value = '  +0.20%  '

$ws.Range("E4").Value = '  -0.35%  '

$ws.Range("E5").Value = '  +7.90%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '247.56'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.21%  '

$ws.Range("E7").Value = '  -0.27%  '

$ws.Range("E8").Value = '  -3.32%  '

$ws.Range("E9").Value = '  +2.43%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '53.00'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +8.16%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0734'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.09%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0990'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.12%  '

$ws.Range("D13").Value = '2.192.14'
$ws.Range("E13").Value = '  +0.11%  '

$ws.Range("E14").Value = '  +2.93%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.719'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.58%  '

$ws.Range("B16").Value = 'Polkadot'
$ws.Range("C16").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.92'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.09%  '

$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '1.913.70'
$ws.Range("E17").Value = '  +1.33%  '

$ws.Range("D18").Value = '35.535.25'
$ws.Range("E18").Value = '  -0.06%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '73.45'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.40%  '

$ws.Range("E20").Value = '  -0.67%  '

$ws.Range("E21").Value = '  +3.49%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '242.87'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.73%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.08'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.95%  '

$ws.Range("E24").Value = '  -0.23%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.33'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.01%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.30'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.85%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '168.39'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.98%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.66'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.69%  '

$ws.Range("E29").Value = '  +2.40%  '

$ws.Range("E30").Value = '  +2.61%  '

$ws.Range("D31").Value = '4.142.54'

$ws.Range("E32").Value = '  +2.28%  '

$ws.Range("E33").Value = '  +1.16%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.93'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +11.83%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.23'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.54%  '

$ws.Range("E36").Value = '  -0.33%  '

$ws.Range("E37").Value = '  -5.81%  '

$ws.Range("E38").Value = '  +10.86%  '

$ws.Range("E39").Value = '  +0.92%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.01'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +15.33%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '99.02'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.59%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.15'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.84%  '

$ws.Range("E43").Value = '  +1.04%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0645'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.92%  '

$ws.Range("D45").Value = '1.351.31'
$ws.Range("E45").Value = '  +0.01%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.46'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.29%  '

$ws.Range("E47").Value = '  +0.34%  '

$ws.Range("E48").Value = '  -0.32%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '45.71'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.22%  '

$ws.Range("E50").Value = '  -0.62%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '12.13'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.48%  '
